# Auto-generated: applies the cryptos-list price/volume refresh described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.553.61"

# Row 3
$ws.Range("D3").Value = "1.565.67"
$ws.Range("E3").Value = "  +0.05%  "

# Row 4
$ws.Range("D4").Value = "'0.989"
$ws.Range("E4").Value = "  -1.70%  "

# Row 5
$ws.Range("D5").Value = "'210.75"
$ws.Range("E5").Value = "  +0.93%  "

# Row 6
$ws.Range("E6").Value = "  +0.11%  "

# Row 7
$ws.Range("E7").Value = "  -1.78%  "

# Row 8
$ws.Range("D8").Value = "'22.74"
$ws.Range("E8").Value = "  +2.89%  "

# Row 9
$ws.Range("E9").Value = "  +0.20%  "

# Row 10
$ws.Range("E10").Value = "  -0.36%  "

# Row 11
$ws.Range("E11").Value = "  +1.08%  "

# Row 12
$ws.Range("D12").Value = "1.789.93"
$ws.Range("E12").Value = "  +0.10%  "

# Row 13
$ws.Range("D13").Value = "1.572.58"
$ws.Range("E13").Value = "  +0.72%  "

# Row 14
$ws.Range("D14").Value = "'3.77"
$ws.Range("E14").Value = "  +0.79%  "

# Row 15
$ws.Range("E15").Value = "  -0.08%  "

# Row 16
$ws.Range("D16").Value = "27.537.42"
$ws.Range("E16").Value = "  +1.81%  "

# Row 17
$ws.Range("D17").Value = "'62.45"
$ws.Range("E17").Value = "  +0.81%  "

# Row 18
$ws.Range("D18").Value = "'225.21"
$ws.Range("E18").Value = "  +4.32%  "

# Row 19
$ws.Range("E19").Value = "  +1.53%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0705"
$ws.Range("E20").Value = "  -0.28%  "

# Row 21
$ws.Range("E21").Value = "  -1.80%  "

# Row 22
$ws.Range("E22").Value = "  -0.81%  "

# Row 23
$ws.Range("D23").Value = "'9.41"
$ws.Range("E23").Value = "  +1.88%  "

# Row 24
$ws.Range("D24").Value = "'1.96"
$ws.Range("E24").Value = "  +0.75%  "

# Row 25
$ws.Range("D25").Value = "'149.69"
$ws.Range("E25").Value = "  -2.58%  "

# Row 26
$ws.Range("D26").Value = "'15.18"
$ws.Range("E26").Value = "  +0.65%  "

# Row 27
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'6.62"
$ws.Range("E27").Value = "  +0.17%  "

# Row 28
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").Value = "'0.108"
$ws.Range("E28").Value = "  +1.97%  "

# Row 29
$ws.Range("D29").Value = "'0.989"
$ws.Range("E29").Value = "  -1.71%  "

# Row 30
$ws.Range("E30").Value = "  +1.14%  "

# Row 31
$ws.Range("E31").Value = "  -0.86%  "

# Row 32
$ws.Range("E32").Value = "  +0.20%  "

# Row 33
$ws.Range("D33").Value = "1.455.24"
$ws.Range("E33").Value = "  +2.01%  "

# Row 34
$ws.Range("D34").Value = "'3.15"
$ws.Range("E34").Value = "  -1.98%  "

# Row 35
$ws.Range("E35").Value = "  +2.25%  "

# Row 36
$ws.Range("D36").Value = "'1.61"
$ws.Range("E36").Value = "  +0.22%  "

# Row 37
$ws.Range("E37").Value = "  -1.29%  "

# Row 38
$ws.Range("E38").Value = "  +0.02%  "

# Row 39
$ws.Range("D39").Value = "'0.541"
$ws.Range("E39").Value = "  +1.47%  "

# Row 40
$ws.Range("D40").Value = "'0.815"
$ws.Range("E40").Value = "  +0.39%  "

# Row 41
$ws.Range("D41").Value = "'5.73"
$ws.Range("E41").Value = "  -1.12%  "

# Row 42
$ws.Range("D42").Value = "'2.34"
$ws.Range("E42").Value = "  +0.55%  "

# Row 43
$ws.Range("E43").Value = "  -1.84%  "

# Row 44
$ws.Range("D44").Value = "'1.85"
$ws.Range("E44").Value = "  +6.15%  "

# Row 45
$ws.Range("D45").Value = "'0.973"
$ws.Range("E45").Value = "  -3.13%  "

# Row 46
$ws.Range("D46").Value = "'64.65"
$ws.Range("E46").Value = "  -0.27%  "

# Row 47
$ws.Range("D47").Value = "1.702.55"
$ws.Range("E47").Value = "  +0.05%  "

# Row 48
$ws.Range("D48").Value = "'86.65"
$ws.Range("E48").Value = "  +0.02%  "

# Row 49
$ws.Range("D49").Value = "'0.0525"
$ws.Range("E49").Value = "  +1.37%  "

# Row 50
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₇0980"
$ws.Range("E50").Value = "  -5.51%  "

# Row 51
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0947"
$ws.Range("E51").Value = "  -1.53%  "

Write-Output "applied cryptos update"
